$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 02:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1527352
$ws.Range("C4").Value = 19579
$ws.Range("D4").Value = 344952
$ws.Range("E4").Value = 1091427

# Row 51 - Chequia
$ws.Range("B51").Value = 8475
$ws.Range("C51").Value = 20
$ws.Range("D51").Value = 5462
$ws.Range("E51").Value = 2715

# Row 53 - Argentina
$ws.Range("B53").Value = 8068
$ws.Range("C53").Value = 263
$ws.Range("E53").Value = 5126
$ws.Range("G53").Value = 10
$ws.Range("H53").Value = 373

# Row 103 - Tunez
$ws.Range("D103").Value = 816
$ws.Range("E103").Value = 176

# Row 126 - Jamaica
$ws.Range("B126").Value = 520
$ws.Range("C126").Value = 9
$ws.Range("D126").Value = 127
$ws.Range("E126").Value = 384

# Row 173 - Malaui
$ws.Range("B173").Value = 70
$ws.Range("C173").Value = 5
$ws.Range("D173").Value = 27
$ws.Range("E173").Value = 40
